{"js": "// Update the date line and every \"NNN\u00d7N=\" expression cell in the\n// multiplication-practice table. Each old value is unique in the\n// document, so a scoped, exact, case-sensitive search/replace per\n// pair is unambiguous and safe.\nconst replacements = [\n  [\"2024-05-28 Tuesday\", \"2024-05-29 Wednesday\"],\n  [\"914\u00d72=\", \"775\u00d76=\"],\n  [\"138\u00d78=\", \"483\u00d72=\"],\n  [\"994\u00d75=\", \"277\u00d77=\"],\n  [\"344\u00d77=\", \"288\u00d78=\"],\n  [\"253\u00d76=\", \"392\u00d77=\"],\n  [\"298\u00d78=\", \"580\u00d72=\"],\n  [\"359\u00d76=\", \"412\u00d75=\"],\n  [\"510\u00d72=\", \"725\u00d73=\"],\n  [\"132\u00d78=\", \"531\u00d75=\"],\n  [\"965\u00d75=\", \"425\u00d74=\"],\n  [\"264\u00d76=\", \"253\u00d75=\"],\n  [\"815\u00d75=\", \"638\u00d72=\"],\n  [\"401\u00d75=\", \"563\u00d73=\"],\n  [\"518\u00d77=\", \"608\u00d76=\"],\n  [\"746\u00d78=\", \"519\u00d76=\"],\n  [\"656\u00d79=\", \"523\u00d74=\"],\n  [\"779\u00d76=\", \"136\u00d73=\"],\n  [\"947\u00d78=\", \"372\u00d72=\"],\n  [\"576\u00d79=\", \"685\u00d77=\"],\n  [\"525\u00d79=\", \"489\u00d73=\"],\n  [\"895\u00d75=\", \"547\u00d77=\"],\n  [\"348\u00d74=\", \"115\u00d77=\"],\n  [\"380\u00d78=\", \"596\u00d72=\"],\n  [\"570\u00d75=\", \"445\u00d79=\"],\n  [\"913\u00d79=\", \"211\u00d74=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const found = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  found.load(\"items\");\n  await context.sync();\n\n  for (const range of found.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Update the date line and every \"NNN\u00d7N=\" expression cell in the\n# multiplication-practice table. Each old value is unique in the\n# document, so a document-wide Find/Replace per pair is unambiguous\n# and safe.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @{ Old = \"2024-05-28 Tuesday\"; New = \"2024-05-29 Wednesday\" },\n    @{ Old = \"914\u00d72=\";             New = \"775\u00d76=\" },\n    @{ Old = \"138\u00d78=\";             New = \"483\u00d72=\" },\n    @{ Old = \"994\u00d75=\";             New = \"277\u00d77=\" },\n    @{ Old = \"344\u00d77=\";             New = \"288\u00d78=\" },\n    @{ Old = \"253\u00d76=\";             New = \"392\u00d77=\" },\n    @{ Old = \"298\u00d78=\";             New = \"580\u00d72=\" },\n    @{ Old = \"359\u00d76=\";             New = \"412\u00d75=\" },\n    @{ Old = \"510\u00d72=\";             New = \"725\u00d73=\" },\n    @{ Old = \"132\u00d78=\";             New = \"531\u00d75=\" },\n    @{ Old = \"965\u00d75=\";             New = \"425\u00d74=\" },\n    @{ Old = \"264\u00d76=\";             New = \"253\u00d75=\" },\n    @{ Old = \"815\u00d75=\";             New = \"638\u00d72=\" },\n    @{ Old = \"401\u00d75=\";             New = \"563\u00d73=\" },\n    @{ Old = \"518\u00d77=\";             New = \"608\u00d76=\" },\n    @{ Old = \"746\u00d78=\";             New = \"519\u00d76=\" },\n    @{ Old = \"656\u00d79=\";             New = \"523\u00d74=\" },\n    @{ Old = \"779\u00d76=\";             New = \"136\u00d73=\" },\n    @{ Old = \"947\u00d78=\";             New = \"372\u00d72=\" },\n    @{ Old = \"576\u00d79=\";             New = \"685\u00d77=\" },\n    @{ Old = \"525\u00d79=\";             New = \"489\u00d73=\" },\n    @{ Old = \"895\u00d75=\";             New = \"547\u00d77=\" },\n    @{ Old = \"348\u00d74=\";             New = \"115\u00d77=\" },\n    @{ Old = \"380\u00d78=\";             New = \"596\u00d72=\" },\n    @{ Old = \"570\u00d75=\";             New = \"445\u00d79=\" },\n    @{ Old = \"913\u00d79=\";             New = \"211\u00d74=\" }\n)\n\nforeach ($pair in $replacements) {\n    $rng = $d.Content\n    $rng.Find.ClearFormatting()\n    $rng.Find.Execute($pair.Old, $false, $true, $false, $false, $false, $true, 1, $false, $pair.New, 2)\n}\n"}
